# Applies the "Add files via upload" edit: adds many new macro/market
# data rows to the three sheets (quarterly / monthly / daily), tweaks a
# couple of existing rows, renames header E1 and widens column A.

$wb = $excel.ActiveWorkbook

$xlLeft = -4131

# ---------------------------------------------------------------------
# Sheet "quarterly" (sheet1)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("quarterly")

# Header row: E1 header text changes from the old Japanese label to "name"
$ws1.Cells.Item(1, 5).Value = "name"

# New data rows 4-10 (style inherited automatically from the column's
# default style on this sheet, so no explicit alignment needed)
$ws1.Cells.Item(4, 1).Value = "A794RX0Q048SBEA"
$ws1.Cells.Item(4, 2).Value = "c_obs"

$ws1.Cells.Item(5, 1).Value = "GPDIC1"
$ws1.Cells.Item(5, 2).Value = "i_obs"

$ws1.Cells.Item(6, 1).Value = "LES1252881600Q"
$ws1.Cells.Item(6, 2).Value = "w_obs"

$ws1.Cells.Item(7, 1).Value = "DPCERD3Q086SBEA"
$ws1.Cells.Item(7, 2).Value = "pi_obs"

$ws1.Cells.Item(8, 1).Value = "BOGZ1FL072052006Q"
$ws1.Cells.Item(8, 2).Value = "r_obs"

$ws1.Cells.Item(9, 1).Value = "GFDEGDQ188S"
$ws1.Cells.Item(9, 2).Value = "b_obs"
$ws1.Cells.Item(9, 6).Value = "対GDP比"

$ws1.Cells.Item(10, 1).Value = "A822RE1Q156NBEA"
$ws1.Cells.Item(10, 2).Value = "g_obs"
$ws1.Cells.Item(10, 6).Value = "対GDP比"

# Column A is widened to fit the longer codes/names
$ws1.Columns.Item(1).ColumnWidth = 28.43

# ---------------------------------------------------------------------
# Sheet "monthly" (sheet2)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("monthly")

$ws2.Cells.Item(1, 5).Value = "name"

# Existing rows just get a renamed "notation" value
$ws2.Cells.Item(2, 2).Value = "epi_obs"
$ws2.Cells.Item(3, 2).Value = "cpi_obs"

# New rows 4-13 - apply the same left/vcenter style used elsewhere on
# the sheet before filling in the values
$ws2.Range("A4:B13").HorizontalAlignment = $xlLeft

$ws2.Cells.Item(4, 1).Value = "CIVPART"
$ws2.Cells.Item(4, 2).Value = "n_obs"

$ws2.Cells.Item(5, 1).Value = "UMCSENT"
$ws2.Cells.Item(5, 2).Value = "sentiment"

$ws2.Cells.Item(6, 1).Value = "M2SL"
$ws2.Cells.Item(6, 2).Value = "money_s"

$ws2.Cells.Item(7, 1).Value = "PAYEMS"
$ws2.Cells.Item(7, 2).Value = "labor_obs"

$ws2.Cells.Item(8, 1).Value = "CSUSHPISA"
$ws2.Cells.Item(8, 2).Value = "hpi"

$ws2.Cells.Item(9, 1).Value = "T10YFFM"
$ws2.Cells.Item(9, 2).Value = "spread_ten"

$ws2.Cells.Item(10, 1).Value = "INTGSBJPM193N"
$ws2.Cells.Item(10, 2).Value = "r_jpn_obs"

$ws2.Cells.Item(11, 1).Value = "JPNWSCNDW01MLSAM"
$ws2.Cells.Item(11, 2).Value = "house_jpn"

$ws2.Cells.Item(12, 1).Value = "MYAGM2JPM189N"
$ws2.Cells.Item(12, 2).Value = "money_s_jpn"

$ws2.Cells.Item(13, 1).Value = "SPASTT01JPM661N"
$ws2.Cells.Item(13, 2).Value = "tosho_jpn"

$ws2.Columns.Item(1).ColumnWidth = 28.43

# ---------------------------------------------------------------------
# Sheet "daily" (sheet3)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("daily")

$ws3.Cells.Item(1, 5).Value = "name"

# Rows 2-4 are re-ordered: wilshire moves up to row 2, sp to row 3,
# oilprice to row 4
$ws3.Cells.Item(2, 1).Value = "WILL5000PRFC"
$ws3.Cells.Item(2, 2).Value = "wilshire"
$ws3.Range("C2:D2").HorizontalAlignment = $xlLeft

$ws3.Cells.Item(3, 1).Value = "SP500"
$ws3.Cells.Item(3, 2).Value = "sp"

$ws3.Cells.Item(4, 1).Value = "DCOILWTICO"
$ws3.Cells.Item(4, 2).Value = "oilprice"

# New rows 5-6
$ws3.Range("A5:B6").HorizontalAlignment = $xlLeft

$ws3.Cells.Item(5, 1).Value = "NASDAQ100"
$ws3.Cells.Item(5, 2).Value = "nasdaq"

$ws3.Cells.Item(6, 1).Value = "NIKKEI225"
$ws3.Cells.Item(6, 2).Value = "nikkei"

$ws3.Columns.Item(1).ColumnWidth = 28.43

# Update the selections on monthly/daily to E2, then re-activate
# quarterly so it stays the tab shown when the workbook is opened.
$ws2.Range("E2").Select()
$ws3.Range("E2").Select()
$ws1.Activate()
